# Updates the crypto price-tracker worksheet with refreshed symbol-list data
# (coinranking scrape), as produced by the "Updated symbol list" GitHub Actions job.
#
# Columns D (Price) and E (Volume(1h)) hold numeric-looking values that are stored
# as plain text in the workbook, so a leading apostrophe is used to stop Excel from
# auto-converting them to numbers/percentages, and the cell style is reset back to
# "Normal" afterwards so no stray number-format is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "E2"; Value = "6.57%"; Numeric = $true },
    @{ Cell = "D3"; Value = "40.21"; Numeric = $true },
    @{ Cell = "E3"; Value = "7.14%"; Numeric = $true },
    @{ Cell = "E4"; Value = "1.87%"; Numeric = $true },
    @{ Cell = "D5"; Value = "0.08092"; Numeric = $true },
    @{ Cell = "E5"; Value = "2.67%"; Numeric = $true },
    @{ Cell = "B6"; Value = "GateToken"; Numeric = $false },
    @{ Cell = "C6"; Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"; Numeric = $false },
    @{ Cell = "D6"; Value = "4.532"; Numeric = $true },
    @{ Cell = "E6"; Value = "2.57%"; Numeric = $true },
    @{ Cell = "B7"; Value = "KuCoinToken"; Numeric = $false },
    @{ Cell = "C7"; Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"; Numeric = $false },
    @{ Cell = "D7"; Value = "8.648"; Numeric = $true },
    @{ Cell = "E7"; Value = "4.55%"; Numeric = $true },
    @{ Cell = "B8"; Value = "FTXToken"; Numeric = $false },
    @{ Cell = "C8"; Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"; Numeric = $false },
    @{ Cell = "D8"; Value = "1.923"; Numeric = $true },
    @{ Cell = "E8"; Value = "0.74%"; Numeric = $true },
    @{ Cell = "B9"; Value = "BTSEToken"; Numeric = $false },
    @{ Cell = "C9"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"; Numeric = $false },
    @{ Cell = "D9"; Value = "2.957"; Numeric = $true },
    @{ Cell = "E9"; Value = "-1.37%"; Numeric = $true },
    @{ Cell = "B10"; Value = "MXToken"; Numeric = $false },
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; Numeric = $false },
    @{ Cell = "D10"; Value = "0.9365"; Numeric = $true },
    @{ Cell = "E10"; Value = "0.00%"; Numeric = $true },
    @{ Cell = "B11"; Value = "LiechtensteinCryptoassetsExchange"; Numeric = $false },
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"; Numeric = $false },
    @{ Cell = "D11"; Value = "0.1340"; Numeric = $true },
    @{ Cell = "E11"; Value = "19.73%"; Numeric = $true },
    @{ Cell = "B12"; Value = "WazirX"; Numeric = $false },
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; Numeric = $false },
    @{ Cell = "D12"; Value = "0.1973"; Numeric = $true },
    @{ Cell = "E12"; Value = "0.03%"; Numeric = $true },
    @{ Cell = "B13"; Value = "MandalaExchangeToken"; Numeric = $false },
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; Numeric = $false },
    @{ Cell = "D13"; Value = "0.09100"; Numeric = $true },
    @{ Cell = "E13"; Value = "0.38%"; Numeric = $true },
    @{ Cell = "B14"; Value = "BitrueCoin"; Numeric = $false },
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; Numeric = $false },
    @{ Cell = "D14"; Value = "0.03505"; Numeric = $true },
    @{ Cell = "E14"; Value = "6.49%"; Numeric = $true },
    @{ Cell = "B15"; Value = "BitMartToken"; Numeric = $false },
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; Numeric = $false },
    @{ Cell = "D15"; Value = "0.09590"; Numeric = $true },
    @{ Cell = "E15"; Value = "-0.02%"; Numeric = $true },
    @{ Cell = "B16"; Value = "BitForexToken"; Numeric = $false },
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; Numeric = $false },
    @{ Cell = "D16"; Value = "0.001374"; Numeric = $true },
    @{ Cell = "E16"; Value = "-0.53%"; Numeric = $true },
    @{ Cell = "B17"; Value = "TigerCash"; Numeric = $false },
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"; Numeric = $false },
    @{ Cell = "D17"; Value = "0.006349"; Numeric = $true },
    @{ Cell = "E17"; Value = "5.00%"; Numeric = $true },
    @{ Cell = "B18"; Value = "LEO"; Numeric = $false },
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; Numeric = $false },
    @{ Cell = "D18"; Value = "3.360"; Numeric = $true },
    @{ Cell = "E18"; Value = "-6.98%"; Numeric = $true },
    @{ Cell = "D19"; Value = "0.3519"; Numeric = $true },
    @{ Cell = "E19"; Value = "3.17%"; Numeric = $true },
    @{ Cell = "D20"; Value = "6.487"; Numeric = $true },
    @{ Cell = "E20"; Value = "1.48%"; Numeric = $true },
    @{ Cell = "D21"; Value = "0.1323"; Numeric = $true },
    @{ Cell = "E21"; Value = "3.25%"; Numeric = $true },
    @{ Cell = "D22"; Value = "0.2579"; Numeric = $true },
    @{ Cell = "E22"; Value = "2.45%"; Numeric = $true },
    @{ Cell = "D23"; Value = "0.04436"; Numeric = $true },
    @{ Cell = "E23"; Value = "0.72%"; Numeric = $true },
    @{ Cell = "D24"; Value = "0.001226"; Numeric = $true },
    @{ Cell = "E24"; Value = "-0.56%"; Numeric = $true },
    @{ Cell = "D25"; Value = "0.004318"; Numeric = $true },
    @{ Cell = "E25"; Value = "-5.80%"; Numeric = $true },
    @{ Cell = "D26"; Value = "0.0001293"; Numeric = $true },
    @{ Cell = "E26"; Value = "-4.96%"; Numeric = $true },
    @{ Cell = "D27"; Value = "0.0003999"; Numeric = $true },
    @{ Cell = "E27"; Value = "0.18%"; Numeric = $true },
    @{ Cell = "D39"; Value = "0.02500"; Numeric = $true },
    @{ Cell = "E39"; Value = "12.98%"; Numeric = $true },
    @{ Cell = "D40"; Value = "0.05191"; Numeric = $true },
    @{ Cell = "E40"; Value = "1.68%"; Numeric = $true },
    @{ Cell = "D41"; Value = "0.007695"; Numeric = $true },
    @{ Cell = "E41"; Value = "3.09%"; Numeric = $true },
    @{ Cell = "D42"; Value = "0.1429"; Numeric = $true },
    @{ Cell = "E42"; Value = "5.69%"; Numeric = $true },
    @{ Cell = "D43"; Value = "0.009224"; Numeric = $true },
    @{ Cell = "E43"; Value = "5.13%"; Numeric = $true },
    @{ Cell = "D44"; Value = "0.002165"; Numeric = $true },
    @{ Cell = "E44"; Value = "1.61%"; Numeric = $true },
    @{ Cell = "D45"; Value = "0.008976"; Numeric = $true },
    @{ Cell = "E45"; Value = "4.09%"; Numeric = $true },
    @{ Cell = "D46"; Value = "0.00006624"; Numeric = $true },
    @{ Cell = "E46"; Value = "1.45%"; Numeric = $true },
    @{ Cell = "D47"; Value = "0.00000000752"; Numeric = $true },
    @{ Cell = "E47"; Value = "0.18%"; Numeric = $true },
    @{ Cell = "D48"; Value = "0.003349"; Numeric = $true },
    @{ Cell = "E48"; Value = "17.13%"; Numeric = $true },
    @{ Cell = "E49"; Value = "148.16%"; Numeric = $true },
    @{ Cell = "D50"; Value = "0.00002105"; Numeric = $true },
    @{ Cell = "E50"; Value = "0.18%"; Numeric = $true },
    @{ Cell = "D51"; Value = "0.0002005"; Numeric = $true },
    @{ Cell = "E51"; Value = "0.18%"; Numeric = $true }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.Numeric) {
        $range.Value = "'" + $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
